# Apply the "kidum" workbook update: swap bitcoin-related keyword/appID rows
# for blockchain / travelpayouts related ones, and clear the last four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (column A = keyword, column B = appID) ---

$ws.Range("A4").Value = "blockchain"
$ws.Range("B4").Value = "block.chain.technology"

$ws.Range("A7").Value = "affiliate marketing"
$ws.Range("B7").Value = "affiliate.marketing.guide"

$ws.Range("A8").Value = "earn passive income"
$ws.Range("B8").Value = "passive.income.nadi.myfirstdrawermenuproject2"

$ws.Range("A9").Value = "Powerful Positive Motivation Quotes"
$ws.Range("B9").Value = "com.sugar.powerfulquotes"

$ws.Range("A10").Value = "blockchain technology"
$ws.Range("B10").Value = "block.chain.technology"

$ws.Range("A11").Value = "affiliate marketing"
$ws.Range("B11").Value = "affiliate.marketing.guide"

$ws.Range("A12").Value = "Powerful Positive Motivation Quotes"
$ws.Range("B12").Value = "com.sugar.powerfulquotes"

$ws.Range("A13").Value = "affiliate marketing"
$ws.Range("B13").Value = "affiliate.marketing.guide"

$ws.Range("A14").Value = "passive income"
$ws.Range("B14").Value = "passive.income.nadi.myfirstdrawermenuproject2"

$ws.Range("A15").Value = "travelpayouts"
$ws.Range("B15").Value = "affiliate.marketing.guide"

$ws.Range("A16").Value = "blockchain"
$ws.Range("B16").Value = "block.chain.technology"

# Rows 17-20 are cleared out entirely in the new layout
$ws.Range("A17:B20").ClearContents()

# --- Row heights that changed along with the content shuffle ---

$ws.Rows.Item(9).RowHeight = 46.5
$ws.Rows.Item(11).RowHeight = 24
$ws.Rows.Item(12).RowHeight = 46.5
$ws.Rows.Item(14).RowHeight = 12.8
$ws.Rows.Item(16).RowHeight = 12.8
$ws.Rows.Item(17).RowHeight = 12.8

# --- Update the view: scroll back to the top and select row 4 ---

$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Rows.Item(4).Select()
